$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 16: B16, C16, E16 become real numbers (were text-like inline strings) ---
$ws.Range("B16").Value = 54446
$ws.Range("C16").Value = 175
$ws.Range("E16").Value = 0

# --- Append new row 17 (still raw/unconverted text values, matching source format) ---
# Use a leading apostrophe to force text interpretation (avoids Excel auto-detecting
# these number/date/percent-looking strings as real numbers), then reset the style
# back to Normal so the quote-prefix formatting doesn't stick to the cell.
$ws.Range("A17").Value = "'2022-01-20"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = "'54446.0"
$ws.Range("B17").Style = "Normal"

$ws.Range("C17").Value = "'706.0"
$ws.Range("C17").Style = "Normal"

$ws.Range("D17").Value = "'1.3%"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "'0"
$ws.Range("E17").Style = "Normal"

$ws.Range("F17").Value = "'"
$ws.Range("F17").Style = "Normal"

$ws.Range("G17").Value = "'"
$ws.Range("G17").Style = "Normal"
